$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.554827
$ws.Range("H2").Value = 10.664481
$ws.Range("I2").Value = 0.2148969460055877
$ws.Range("J2").Value = 0.2148969460055877
$ws.Range("M2").Value = 108.845309
$ws.Range("N2").Value = 326.535927
$ws.Range("O2").Value = 0.3930063530400584
$ws.Range("P2").Value = 0.3930063530400583
$ws.Range("Q2").Value = 386.9262432565429
$ws.Range("R2").Value = 3482.336189308887
$ws.Range("S2").Value = 0.08445586502910236
$ws.Range("T2").Value = 0.08445586502910234
$ws.Range("G3").Value = 3.554827
$ws.Range("H3").Value = 10.664481
$ws.Range("I3").Value = 0.2148969460055877
$ws.Range("J3").Value = 0.2148969460055877
$ws.Range("O3").Value = 0.4513455929560021
$ws.Range("P3").Value = 0.451345592956002
$ws.Range("Q3").Value = 444.3629303749759
$ws.Range("R3").Value = 3999.266373374783
$ws.Range("S3").Value = 0.09699278951932594
$ws.Range("T3").Value = 0.09699278951932591
$ws.Range("G4").Value = 3.554827
$ws.Range("H4").Value = 10.664481
$ws.Range("I4").Value = 0.2148969460055877
$ws.Range("J4").Value = 0.2148969460055877
$ws.Range("M4").Value = 18.88973866666667
$ws.Range("N4").Value = 56.669216
$ws.Range("O4").Value = 0.06820493571538706
$ws.Range("P4").Value = 0.06820493571538705
$ws.Range("Q4").Value = 67.14975303521065
$ws.Range("R4").Value = 604.3477773168959
$ws.Range("S4").Value = 0.01465703238774411
$ws.Range("T4").Value = 0.01465703238774411
$ws.Range("G5").Value = 3.554827
$ws.Range("H5").Value = 10.664481
$ws.Range("I5").Value = 0.2148969460055877
$ws.Range("J5").Value = 0.2148969460055877
$ws.Range("M5").Value = 13.37928666666667
$ws.Range("N5").Value = 40.13786
$ws.Range("O5").Value = 0.04830841776694433
$ws.Range("P5").Value = 0.04830841776694433
$ws.Range("Q5").Value = 47.56104948340666
$ws.Range("R5").Value = 428.04944535066
$ws.Range("S5").Value = 0.01038133144447841
$ws.Range("T5").Value = 0.01038133144447841
$ws.Range("G6").Value = 3.554827
$ws.Range("H6").Value = 10.664481
$ws.Range("I6").Value = 0.2148969460055877
$ws.Range("J6").Value = 0.2148969460055877
$ws.Range("M6").Value = 10.83857433333333
$ws.Range("N6").Value = 32.515723
$ws.Range("O6").Value = 0.03913470052160829
$ws.Range("P6").Value = 0.03913470052160829
$ws.Range("Q6").Value = 38.52925668164033
$ws.Range("R6").Value = 346.7633101347629
$ws.Range("S6").Value = 0.008409927624936902
$ws.Range("T6").Value = 0.0084099276249369
$ws.Range("I7").Value = 0.3107709374420163
$ws.Range("J7").Value = 0.3107709374420163
$ws.Range("M7").Value = 108.845309
$ws.Range("N7").Value = 326.535927
$ws.Range("O7").Value = 0.3930063530400584
$ws.Range("P7").Value = 0.3930063530400583
$ws.Range("Q7").Value = 559.5492796562446
$ws.Range("R7").Value = 5035.943516906202
$ws.Range("S7").Value = 0.1221349527549269
$ws.Range("T7").Value = 0.1221349527549269
$ws.Range("I8").Value = 0.3107709374420163
$ws.Range("J8").Value = 0.3107709374420163
$ws.Range("O8").Value = 0.4513455929560021
$ws.Range("P8").Value = 0.451345592956002
$ws.Range("S8").Value = 0.1402650930332595
$ws.Range("T8").Value = 0.1402650930332594
$ws.Range("I9").Value = 0.3107709374420163
$ws.Range("J9").Value = 0.3107709374420163
$ws.Range("M9").Value = 18.88973866666667
$ws.Range("N9").Value = 56.669216
$ws.Range("O9").Value = 0.06820493571538706
$ws.Range("P9").Value = 0.06820493571538705
$ws.Range("Q9").Value = 97.10790259071288
$ws.Range("R9").Value = 873.9711233164159
$ws.Range("S9").Value = 0.0211961118104433
$ws.Range("T9").Value = 0.02119611181044329
$ws.Range("I10").Value = 0.3107709374420163
$ws.Range("J10").Value = 0.3107709374420163
$ws.Range("M10").Value = 13.37928666666667
$ws.Range("N10").Value = 40.13786
$ws.Range("O10").Value = 0.04830841776694433
$ws.Range("P10").Value = 0.04830841776694433
$ws.Range("Q10").Value = 68.77990687359555
$ws.Range("R10").Value = 619.01916186236
$ws.Range("S10").Value = 0.01501285227577384
$ws.Range("T10").Value = 0.01501285227577384
$ws.Range("I11").Value = 0.3107709374420163
$ws.Range("J11").Value = 0.3107709374420163
$ws.Range("M11").Value = 10.83857433333333
$ws.Range("N11").Value = 32.515723
$ws.Range("O11").Value = 0.03913470052160829
$ws.Range("P11").Value = 0.03913470052160829
$ws.Range("Q11").Value = 55.71867558129978
$ws.Range("R11").Value = 501.468080231698
$ws.Range("S11").Value = 0.01216192756761277
$ws.Range("T11").Value = 0.01216192756761277
$ws.Range("G12").Value = 2.615693
$ws.Range("H12").Value = 7.847079000000001
$ws.Range("I12").Value = 0.1581242736673807
$ws.Range("J12").Value = 0.1581242736673807
$ws.Range("M12").Value = 108.845309
$ws.Range("N12").Value = 326.535927
$ws.Range("O12").Value = 0.3930063530400584
$ws.Range("P12").Value = 0.3930063530400583
$ws.Range("Q12").Value = 284.7059128341371
$ws.Range("R12").Value = 2562.353215507233
$ws.Range("S12").Value = 0.06214384412112541
$ws.Range("T12").Value = 0.0621438441211254
$ws.Range("G13").Value = 2.615693
$ws.Range("H13").Value = 7.847079000000001
$ws.Range("I13").Value = 0.1581242736673807
$ws.Range("J13").Value = 0.1581242736673807
$ws.Range("O13").Value = 0.4513455929560021
$ws.Range("P13").Value = 0.451345592956002
$ws.Range("Q13").Value = 326.968655982784
$ws.Range("R13").Value = 2942.717903845056
$ws.Range("S13").Value = 0.07136869405914106
$ws.Range("T13").Value = 0.07136869405914105
$ws.Range("G14").Value = 2.615693
$ws.Range("H14").Value = 7.847079000000001
$ws.Range("I14").Value = 0.1581242736673807
$ws.Range("J14").Value = 0.1581242736673807
$ws.Range("M14").Value = 18.88973866666667
$ws.Range("N14").Value = 56.669216
$ws.Range("O14").Value = 0.06820493571538706
$ws.Range("P14").Value = 0.06820493571538705
$ws.Range("Q14").Value = 49.40975720222934
$ws.Range("R14").Value = 444.687814820064
$ws.Range("S14").Value = 0.01078485592052597
$ws.Range("T14").Value = 0.01078485592052597
$ws.Range("G15").Value = 2.615693
$ws.Range("H15").Value = 7.847079000000001
$ws.Range("I15").Value = 0.1581242736673807
$ws.Range("J15").Value = 0.1581242736673807
$ws.Range("M15").Value = 13.37928666666667
$ws.Range("N15").Value = 40.13786
$ws.Range("O15").Value = 0.04830841776694433
$ws.Range("P15").Value = 0.04830841776694433
$ws.Range("Q15").Value = 34.99610647899334
$ws.Range("R15").Value = 314.9649583109401
$ws.Range("S15").Value = 0.00763873347141846
$ws.Range("T15").Value = 0.007638733471418459
$ws.Range("G16").Value = 2.615693
$ws.Range("H16").Value = 7.847079000000001
$ws.Range("I16").Value = 0.1581242736673807
$ws.Range("J16").Value = 0.1581242736673807
$ws.Range("M16").Value = 10.83857433333333
$ws.Range("N16").Value = 32.515723
$ws.Range("O16").Value = 0.03913470052160829
$ws.Range("P16").Value = 0.03913470052160829
$ws.Range("Q16").Value = 28.35038301367967
$ws.Range("R16").Value = 255.153447123117
$ws.Range("S16").Value = 0.006188146095169774
$ws.Range("T16").Value = 0.006188146095169773
$ws.Range("G17").Value = 4.248598333333334
$ws.Range("H17").Value = 12.745795
$ws.Range("I17").Value = 0.2568369168563656
$ws.Range("J17").Value = 0.2568369168563656
$ws.Range("M17").Value = 108.845309
$ws.Range("N17").Value = 326.535927
$ws.Range("O17").Value = 0.3930063530400584
$ws.Range("P17").Value = 0.3930063530400583
$ws.Range("Q17").Value = 462.4399984085517
$ws.Range("R17").Value = 4161.959985676966
$ws.Range("S17").Value = 0.1009385400197729
$ws.Range("T17").Value = 0.1009385400197729
$ws.Range("G18").Value = 4.248598333333334
$ws.Range("H18").Value = 12.745795
$ws.Range("I18").Value = 0.2568369168563656
$ws.Range("J18").Value = 0.2568369168563656
$ws.Range("O18").Value = 0.4513455929560021
$ws.Range("P18").Value = 0.451345592956002
$ws.Range("Q18").Value = 531.0862118989867
$ws.Range("R18").Value = 4779.775907090881
$ws.Range("S18").Value = 0.1159222105315277
$ws.Range("T18").Value = 0.1159222105315277
$ws.Range("G19").Value = 4.248598333333334
$ws.Range("H19").Value = 12.745795
$ws.Range("I19").Value = 0.2568369168563656
$ws.Range("J19").Value = 0.2568369168563656
$ws.Range("M19").Value = 18.88973866666667
$ws.Range("N19").Value = 56.669216
$ws.Range("O19").Value = 0.06820493571538706
$ws.Range("P19").Value = 0.06820493571538705
$ws.Range("Q19").Value = 80.25491221630223
$ws.Range("R19").Value = 722.29420994672
$ws.Range("S19").Value = 0.01751754540352662
$ws.Range("T19").Value = 0.01751754540352662
$ws.Range("G20").Value = 4.248598333333334
$ws.Range("H20").Value = 12.745795
$ws.Range("I20").Value = 0.2568369168563656
$ws.Range("J20").Value = 0.2568369168563656
$ws.Range("M20").Value = 13.37928666666667
$ws.Range("N20").Value = 40.13786
$ws.Range("O20").Value = 0.04830841776694433
$ws.Range("P20").Value = 0.04830841776694433
$ws.Range("Q20").Value = 56.8432150331889
$ws.Range("R20").Value = 511.5889352987001
$ws.Range("S20").Value = 0.01240738507747126
$ws.Range("T20").Value = 0.01240738507747125
$ws.Range("G21").Value = 4.248598333333334
$ws.Range("H21").Value = 12.745795
$ws.Range("I21").Value = 0.2568369168563656
$ws.Range("J21").Value = 0.2568369168563656
$ws.Range("M21").Value = 10.83857433333333
$ws.Range("N21").Value = 32.515723
$ws.Range("O21").Value = 0.03913470052160829
$ws.Range("P21").Value = 0.03913470052160829
$ws.Range("Q21").Value = 46.04874884830945
$ws.Range("R21").Value = 414.4387396347851
$ws.Range("S21").Value = 0.01005123582406708
$ws.Range("T21").Value = 0.01005123582406707
$ws.Range("G22").Value = 0.9821143333333332
$ws.Range("H22").Value = 2.946343
$ws.Range("I22").Value = 0.0593709260286498
$ws.Range("J22").Value = 0.0593709260286498
$ws.Range("M22").Value = 108.845309
$ws.Range("N22").Value = 326.535927
$ws.Range("O22").Value = 0.3930063530400584
$ws.Range("P22").Value = 0.3930063530400583
$ws.Range("Q22").Value = 106.8985380849956
$ws.Range("R22").Value = 962.0868427649609
$ws.Range("S22").Value = 0.02333315111513073
$ws.Range("T22").Value = 0.02333315111513073
$ws.Range("G23").Value = 0.9821143333333332
$ws.Range("H23").Value = 2.946343
$ws.Range("I23").Value = 0.0593709260286498
$ws.Range("J23").Value = 0.0593709260286498
$ws.Range("O23").Value = 0.4513455929560021
$ws.Range("P23").Value = 0.451345592956002
$ws.Range("Q23").Value = 122.7669315899946
$ws.Range("R23").Value = 1104.902384309952
$ws.Range("S23").Value = 0.02679680581274788
$ws.Range("T23").Value = 0.02679680581274788
$ws.Range("G24").Value = 0.9821143333333332
$ws.Range("H24").Value = 2.946343
$ws.Range("I24").Value = 0.0593709260286498
$ws.Range("J24").Value = 0.0593709260286498
$ws.Range("M24").Value = 18.88973866666667
$ws.Range("N24").Value = 56.669216
$ws.Range("O24").Value = 0.06820493571538706
$ws.Range("P24").Value = 0.06820493571538705
$ws.Range("Q24").Value = 18.55188309745422
$ws.Range("R24").Value = 166.966947877088
$ws.Range("S24").Value = 0.00404939019314706
$ws.Range("T24").Value = 0.00404939019314706
$ws.Range("G25").Value = 0.9821143333333332
$ws.Range("H25").Value = 2.946343
$ws.Range("I25").Value = 0.0593709260286498
$ws.Range("J25").Value = 0.0593709260286498
$ws.Range("M25").Value = 13.37928666666667
$ws.Range("N25").Value = 40.13786
$ws.Range("O25").Value = 0.04830841776694433
$ws.Range("P25").Value = 0.04830841776694433
$ws.Range("Q25").Value = 13.13998920510889
$ws.Range("R25").Value = 118.25990284598
$ws.Range("S25").Value = 0.002868115497802364
$ws.Range("T25").Value = 0.002868115497802364
$ws.Range("G26").Value = 0.9821143333333332
$ws.Range("H26").Value = 2.946343
$ws.Range("I26").Value = 0.0593709260286498
$ws.Range("J26").Value = 0.0593709260286498
$ws.Range("M26").Value = 10.83857433333333
$ws.Range("N26").Value = 32.515723
$ws.Range("O26").Value = 0.03913470052160829
$ws.Range("P26").Value = 0.03913470052160829
$ws.Range("Q26").Value = 10.64471920566544
$ws.Range("R26").Value = 95.80247285098899
$ws.Range("S26").Value = 0.002323463409821769
$ws.Range("T26").Value = 0.002323463409821769
